$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = 10.13
    4 = 12.28
    5 = 15.72
    6 = 2.8
    7 = 8.59
    8 = 7.56
    9 = 5.7
    10 = 24.16
    11 = 1.27
    12 = 2.82
    13 = 1.2
    14 = 6.45
    16 = 4.41
    17 = 16.71
    18 = 3.42
    19 = 12.64
    20 = 4.81
    21 = 1.51
    22 = 6.69
    23 = 2.9
    24 = 3.61
    26 = 4.91
    29 = 0.58
    30 = 0.93
    31 = 2
    34 = 9.960000000000001
    36 = 1.85
    37 = 5.22
    38 = 25.33
    39 = 0.85
    41 = 2.05
    42 = 3.64
    43 = 3.16
    44 = 0.76
    46 = 3.62
    47 = 8.140000000000001
    53 = 3.47
    54 = 3.53
    56 = 11.01
    57 = 8.51
    58 = 16.45
    60 = 8.74
    61 = 8.970000000000001
    62 = 3.5
    63 = 1.66
    64 = 5.58
    65 = 3.81
    66 = 0.97
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item([int]$row, 4).Value = $updates[$row]
}

